$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New row 41 inherits the same formatting as the row above it (row 40):
# copy formats first so borders/fill/number-format line up, then fix up
# the value cells individually.
$ws.Range("A40:E40").Copy()
$ws.Range("A41:E41").PasteSpecial(-4122)

$ws.Range("A41").Value = "IAM040"

# Column B on this new row uses the plain (no border/fill) wrap-text style
# that the B column itself defaults to, instead of row 40's bordered style.
$ws.Range("B41").Borders.LineStyle = -4142
$ws.Range("B41").Interior.Pattern = -4142
$ws.Range("B41").WrapText = $true
$ws.Range("B41").Value = "OPQA-5372||OPQA-5373||OPQA-4252"

$ws.Range("C41").Value = "Verify that error message " + [char]34 + "Please enter an email address." + [char]34 + " should be displayed in red color when user not enter email address in email text field for neon login page.||Verify that error message " + [char]34 + "Please enter a password." + [char]34 + " should be displayed in red color when user not enter email address in email text field for neon login page||Verify that error message " + [char]34 + " Please enter a valid email address." + [char]34 + " should be displayed in red color when user enters email address in wrong format"

$ws.Range("D41").Value = "Y"
$ws.Range("E41").Value = "PASS"

$ws.Rows.Item(41).RowHeight = 86.4

$ws.Range("C41").Select()
